# Auto-generated edit script applying numeric (and one status-text) corrections
# to the F ("想去人数") and G ("最低票价") columns across all four sheets,
# matching the authoritative diff between before.xlsx and the target workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1904
$ws.Range("F3").Value = 29
$ws.Range("F4").Value = 79
$ws.Range("F5").Value = 799
$ws.Range("F8").Value = 981
$ws.Range("F9").Value = 1657
$ws.Range("F10").Value = 1306
$ws.Range("F11").Value = 1595
$ws.Range("F13").Value = 1598
$ws.Range("F14").Value = 360
$ws.Range("F15").Value = 1728
$ws.Range("F17").Value = 1162
$ws.Range("F19").Value = 2018
$ws.Range("F20").Value = 277
$ws.Range("F21").Value = 836
$ws.Range("F22").Value = 1023
$ws.Range("F24").Value = 5
$ws.Range("F25").Value = 1325
$ws.Range("F26").Value = 1111
$ws.Range("F27").Value = 101
$ws.Range("F29").Value = 1234
$ws.Range("F30").Value = 921
$ws.Range("F31").Value = 1220
$ws.Range("F32").Value = 60
$ws.Range("F33").Value = 1166
$ws.Range("F34").Value = 320
$ws.Range("F35").Value = 92
$ws.Range("F36").Value = 908
$ws.Range("F38").Value = 1739
$ws.Range("F39").Value = 394
$ws.Range("F41").Value = 132
$ws.Range("F42").Value = 2112
$ws.Range("F43").Value = 107
$ws.Range("F44").Value = 851
$ws.Range("F45").Value = 732
$ws.Range("F47").Value = 814
$ws.Range("F48").Value = 125

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 2633
$ws.Range("F10").Value = 424
$ws.Range("F11").Value = 736
$ws.Range("G11").Value = "不可售"
$ws.Range("F15").Value = 22
$ws.Range("F27").Value = 198
$ws.Range("F28").Value = 261
$ws.Range("F32").Value = 61
$ws.Range("F34").Value = 34
$ws.Range("F36").Value = 19
$ws.Range("F41").Value = 73

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2956
$ws.Range("F6").Value = 4732
$ws.Range("F7").Value = 159
$ws.Range("F10").Value = 821
$ws.Range("F11").Value = 491
$ws.Range("F12").Value = 476
$ws.Range("F13").Value = 1198
$ws.Range("F14").Value = 333
$ws.Range("F15").Value = 852

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1904
$ws.Range("F4").Value = 29
$ws.Range("F5").Value = 4732
$ws.Range("F7").Value = 821
$ws.Range("F8").Value = 491
$ws.Range("F9").Value = 476
$ws.Range("F10").Value = 476
$ws.Range("F11").Value = 1198
$ws.Range("F13").Value = 981
$ws.Range("F14").Value = 1657
$ws.Range("F15").Value = 1306
$ws.Range("F16").Value = 1595
$ws.Range("F18").Value = 1598
$ws.Range("F21").Value = 1728
$ws.Range("F22").Value = 1162
$ws.Range("F23").Value = 852
$ws.Range("F24").Value = 852
$ws.Range("F25").Value = 2018
$ws.Range("F26").Value = 277
$ws.Range("F27").Value = 836
$ws.Range("F28").Value = 1023
$ws.Range("F30").Value = 1325
$ws.Range("F32").Value = 1111
$ws.Range("F33").Value = 101
$ws.Range("F34").Value = 1234
$ws.Range("F35").Value = 921
$ws.Range("F36").Value = 1220
$ws.Range("F37").Value = 60
$ws.Range("F40").Value = 1166
$ws.Range("F41").Value = 320
$ws.Range("F42").Value = 908
$ws.Range("F44").Value = 1739
$ws.Range("F46").Value = 2112
$ws.Range("F47").Value = 107
$ws.Range("F48").Value = 851
$ws.Range("F49").Value = 732
$ws.Range("F50").Value = 814
$ws.Range("F51").Value = 125

